# Auto update Excel log: append newly logged sensor readings for 2026-02-06
$wb = $excel.ActiveWorkbook

# ---- PIR sheet ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("A434:A446").NumberFormat = "@"

$ws.Range("A434").Value = '2026-02-06'
$ws.Range("B434").Value = '10:15:56'
$ws.Range("C434").Value = '10:00'
$ws.Range("D434").Value = 'Bathroom'
$ws.Range("E434").Value = 'No Motion'
$ws.Range("F434").Value = 'Inactive'

$ws.Range("A435").Value = '2026-02-06'
$ws.Range("B435").Value = '10:16:01'
$ws.Range("C435").Value = '10:00'
$ws.Range("D435").Value = 'Bathroom'
$ws.Range("E435").Value = 'No Motion'
$ws.Range("F435").Value = 'Inactive'

$ws.Range("A436").Value = '2026-02-06'
$ws.Range("B436").Value = '10:16:04'
$ws.Range("C436").Value = '10:00'
$ws.Range("D436").Value = 'Bathroom'
$ws.Range("E436").Value = 'No Motion'
$ws.Range("F436").Value = 'Inactive'

$ws.Range("A437").Value = '2026-02-06'
$ws.Range("B437").Value = '10:16:11'
$ws.Range("C437").Value = '10:00'
$ws.Range("D437").Value = 'Bathroom'
$ws.Range("E437").Value = 'No Motion'
$ws.Range("F437").Value = 'Inactive'

$ws.Range("A438").Value = '2026-02-06'
$ws.Range("B438").Value = '10:16:16'
$ws.Range("C438").Value = '10:00'
$ws.Range("D438").Value = 'Bathroom'
$ws.Range("E438").Value = 'No Motion'
$ws.Range("F438").Value = 'Inactive'

$ws.Range("A439").Value = '2026-02-06'
$ws.Range("B439").Value = '10:16:21'
$ws.Range("C439").Value = '10:00'
$ws.Range("D439").Value = 'Bathroom'
$ws.Range("E439").Value = 'No Motion'
$ws.Range("F439").Value = 'Inactive'

$ws.Range("A440").Value = '2026-02-06'
$ws.Range("B440").Value = '10:16:26'
$ws.Range("C440").Value = '10:00'
$ws.Range("D440").Value = 'Bathroom'
$ws.Range("E440").Value = 'No Motion'
$ws.Range("F440").Value = 'Inactive'

$ws.Range("A441").Value = '2026-02-06'
$ws.Range("B441").Value = '10:16:27'
$ws.Range("C441").Value = '10:00'
$ws.Range("D441").Value = 'Bathroom'
$ws.Range("E441").Value = 'Motion Detected'
$ws.Range("F441").Value = 'Active'

$ws.Range("A442").Value = '2026-02-06'
$ws.Range("B442").Value = '10:16:32'
$ws.Range("C442").Value = '10:00'
$ws.Range("D442").Value = 'Bathroom'
$ws.Range("E442").Value = 'No Motion'
$ws.Range("F442").Value = 'Inactive'

$ws.Range("A443").Value = '2026-02-06'
$ws.Range("B443").Value = '10:16:38'
$ws.Range("C443").Value = '10:00'
$ws.Range("D443").Value = 'Bathroom'
$ws.Range("E443").Value = 'No Motion'
$ws.Range("F443").Value = 'Inactive'

$ws.Range("A444").Value = '2026-02-06'
$ws.Range("B444").Value = '10:16:43'
$ws.Range("C444").Value = '10:00'
$ws.Range("D444").Value = 'Bathroom'
$ws.Range("E444").Value = 'No Motion'
$ws.Range("F444").Value = 'Inactive'

$ws.Range("A445").Value = '2026-02-06'
$ws.Range("B445").Value = '10:16:46'
$ws.Range("C445").Value = '10:00'
$ws.Range("D445").Value = 'Bathroom'
$ws.Range("E445").Value = 'Motion Detected'
$ws.Range("F445").Value = 'Active'

$ws.Range("A446").Value = '2026-02-06'
$ws.Range("B446").Value = '10:16:54'
$ws.Range("C446").Value = '10:00'
$ws.Range("D446").Value = 'Bathroom'
$ws.Range("E446").Value = 'No Motion'
$ws.Range("F446").Value = 'Inactive'

# ---- Humidity sheet ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("A298:A308").NumberFormat = "@"
$ws.Range("E298:E308").NumberFormat = "@"

$ws.Range("A298").Value = '2026-02-06'
$ws.Range("B298").Value = '10:15:59'
$ws.Range("C298").Value = '10:00'
$ws.Range("D298").Value = 'Bathroom'
$ws.Range("E298").Value = '68.8%'
$ws.Range("F298").Value = 'Active'

$ws.Range("A299").Value = '2026-02-06'
$ws.Range("B299").Value = '10:16:09'
$ws.Range("C299").Value = '10:00'
$ws.Range("D299").Value = 'Bathroom'
$ws.Range("E299").Value = '67.7%'
$ws.Range("F299").Value = 'Active'

$ws.Range("A300").Value = '2026-02-06'
$ws.Range("B300").Value = '10:16:14'
$ws.Range("C300").Value = '10:00'
$ws.Range("D300").Value = 'Bathroom'
$ws.Range("E300").Value = '68.7%'
$ws.Range("F300").Value = 'Active'

$ws.Range("A301").Value = '2026-02-06'
$ws.Range("B301").Value = '10:16:19'
$ws.Range("C301").Value = '10:00'
$ws.Range("D301").Value = 'Bathroom'
$ws.Range("E301").Value = '67.8%'
$ws.Range("F301").Value = 'Active'

$ws.Range("A302").Value = '2026-02-06'
$ws.Range("B302").Value = '10:16:24'
$ws.Range("C302").Value = '10:00'
$ws.Range("D302").Value = 'Bathroom'
$ws.Range("E302").Value = '68.9%'
$ws.Range("F302").Value = 'Active'

$ws.Range("A303").Value = '2026-02-06'
$ws.Range("B303").Value = '10:16:29'
$ws.Range("C303").Value = '10:00'
$ws.Range("D303").Value = 'Bathroom'
$ws.Range("E303").Value = '68.0%'
$ws.Range("F303").Value = 'Active'

$ws.Range("A304").Value = '2026-02-06'
$ws.Range("B304").Value = '10:16:34'
$ws.Range("C304").Value = '10:00'
$ws.Range("D304").Value = 'Bathroom'
$ws.Range("E304").Value = '68.9%'
$ws.Range("F304").Value = 'Active'

$ws.Range("A305").Value = '2026-02-06'
$ws.Range("B305").Value = '10:16:39'
$ws.Range("C305").Value = '10:00'
$ws.Range("D305").Value = 'Bathroom'
$ws.Range("E305").Value = '67.8%'
$ws.Range("F305").Value = 'Active'

$ws.Range("A306").Value = '2026-02-06'
$ws.Range("B306").Value = '10:16:44'
$ws.Range("C306").Value = '10:00'
$ws.Range("D306").Value = 'Bathroom'
$ws.Range("E306").Value = '68.7%'
$ws.Range("F306").Value = 'Active'

$ws.Range("A307").Value = '2026-02-06'
$ws.Range("B307").Value = '10:16:49'
$ws.Range("C307").Value = '10:00'
$ws.Range("D307").Value = 'Bathroom'
$ws.Range("E307").Value = '67.7%'
$ws.Range("F307").Value = 'Active'

$ws.Range("A308").Value = '2026-02-06'
$ws.Range("B308").Value = '10:16:55'
$ws.Range("C308").Value = '10:00'
$ws.Range("D308").Value = 'Bathroom'
$ws.Range("E308").Value = '68.8%'
$ws.Range("F308").Value = 'Active'

# ---- Temperature sheet ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("A298:A307").NumberFormat = "@"

$ws.Range("A298").Value = '2026-02-06'
$ws.Range("B298").Value = '10:16:00'
$ws.Range("C298").Value = '10:00'
$ws.Range("D298").Value = 'Bathroom'
$ws.Range("E298").Value = '28.2C'
$ws.Range("F298").Value = 'Active'

$ws.Range("A299").Value = '2026-02-06'
$ws.Range("B299").Value = '10:16:10'
$ws.Range("C299").Value = '10:00'
$ws.Range("D299").Value = 'Bathroom'
$ws.Range("E299").Value = '28.1C'
$ws.Range("F299").Value = 'Active'

$ws.Range("A300").Value = '2026-02-06'
$ws.Range("B300").Value = '10:16:15'
$ws.Range("C300").Value = '10:00'
$ws.Range("D300").Value = 'Bathroom'
$ws.Range("E300").Value = '28.1C'
$ws.Range("F300").Value = 'Active'

$ws.Range("A301").Value = '2026-02-06'
$ws.Range("B301").Value = '10:16:20'
$ws.Range("C301").Value = '10:00'
$ws.Range("D301").Value = 'Bathroom'
$ws.Range("E301").Value = '28.2C'
$ws.Range("F301").Value = 'Active'

$ws.Range("A302").Value = '2026-02-06'
$ws.Range("B302").Value = '10:16:25'
$ws.Range("C302").Value = '10:00'
$ws.Range("D302").Value = 'Bathroom'
$ws.Range("E302").Value = '28.1C'
$ws.Range("F302").Value = 'Active'

$ws.Range("A303").Value = '2026-02-06'
$ws.Range("B303").Value = '10:16:30'
$ws.Range("C303").Value = '10:00'
$ws.Range("D303").Value = 'Bathroom'
$ws.Range("E303").Value = '28.2C'
$ws.Range("F303").Value = 'Active'

$ws.Range("A304").Value = '2026-02-06'
$ws.Range("B304").Value = '10:16:35'
$ws.Range("C304").Value = '10:00'
$ws.Range("D304").Value = 'Bathroom'
$ws.Range("E304").Value = '28.3C'
$ws.Range("F304").Value = 'Active'

$ws.Range("A305").Value = '2026-02-06'
$ws.Range("B305").Value = '10:16:40'
$ws.Range("C305").Value = '10:00'
$ws.Range("D305").Value = 'Bathroom'
$ws.Range("E305").Value = '28.2C'
$ws.Range("F305").Value = 'Active'

$ws.Range("A306").Value = '2026-02-06'
$ws.Range("B306").Value = '10:16:45'
$ws.Range("C306").Value = '10:00'
$ws.Range("D306").Value = 'Bathroom'
$ws.Range("E306").Value = '28.2C'
$ws.Range("F306").Value = 'Active'

$ws.Range("A307").Value = '2026-02-06'
$ws.Range("B307").Value = '10:16:50'
$ws.Range("C307").Value = '10:00'
$ws.Range("D307").Value = 'Bathroom'
$ws.Range("E307").Value = '28.1C'
$ws.Range("F307").Value = 'Active'

